$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data row (row 2) with new test data values
$ws.Range("B2").Value = "Fernando Alonso"
$ws.Range("A2").Value = "TEST-PLATE"
$ws.Range("C2").Value = 2001
$ws.Range("D2").Value = "Green"

# Update the active selection to D3 (matches sheetView selection in diff)
$ws.Range("D3").Select()
